$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.411.36'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '1.848.73'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9997'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '240.50'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6281'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07690'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.2921'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.24%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '25.02'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.83%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07746'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.45%  '
$ws.Range('D12').Value = '1.870.37'
$ws.Range('E12').Value = '  +1.80%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.035'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.67%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.00001085'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.83%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6824'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.36%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '83.61'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.38%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '6.194'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('D18').Value = '29.442.64'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.40'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.0000'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.462'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '157.60'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.46%  '
$ws.Range('E25').Value = '  -1.03%  '
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '17.72'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.74%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.352'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.97%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.462'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05632'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.047'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.846'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('E34').Value = '  +0.61%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7014'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.27%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.593'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('D37').Value = '1.225.97'
$ws.Range('E37').Value = '  -1.67%  '
$ws.Range('E38').Value = '  -1.00%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.751'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.447'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.94%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9054'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.37%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.000'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '2.026.45'
$ws.Range('E43').Value = '  +1.46%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '101.85'
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '66.09'
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.188'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.34%  '
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.4021'
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.1156'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.06%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.008'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.95%  '
$ws.Range('E51').Value = '  +0.33%  '
